# The "(calculated)" label in A6 was removed from the LFO-Analog rates sheet,
# which also drops that now-unused entry from the shared string table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A6").ClearContents()

# Reset the saved selection back to the top-left cell instead of the old C7.
$ws.Range("A1").Select()
